$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.581.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.663.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.06%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.94%  '
$ws.Range("E6").Value = '  -2.78%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.32'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("E10").Value = '  -2.70%  '
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.900.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.666.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.95%  '
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.568'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.576.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '240.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0731'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.28%  '
$ws.Range("E20").Value = '  -3.59%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("E23").Value = '  -3.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.98%  '
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  -2.63%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0504'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  -2.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.460.62'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.93%  '
$ws.Range("E34").Value = '  -4.47%  '
$ws.Range("E35").Value = '  -4.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.931'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.90%  '
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("E39").Value = '  -5.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '70.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E41").Value = '  -5.05%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  -4.39%  '
$ws.Range("E44").Value = '  -3.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.795'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.807.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.02%  '
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("E49").Value = '  -5.03%  '
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.65%  '
